$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valores")

$ws.Range("I1").Value = 6

$ws.Range("G2").Value = 1
$ws.Range("J2").Value = 4
$ws.Range("L2").Value = 5

$ws.Range("I4").Value = 1

$ws.Range("I6").Value = 3

$ws.Range("J10").Value = 1

$ws.Range("J15").Value = 1
$ws.Range("L15").Value = 5

$ws.Range("G16").Value = 1
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 6
$ws.Range("L16").Value = 10
